$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from H1 (bold, bordered, centered header style)
# onto the two new header cells I1 and J1.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Set the new header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New column data (I = I0, J = IF) for rows 2-11
$values = @{
    2  = @(8, 8)
    3  = @(9, 9)
    4  = @(8, 9)
    5  = @(8, 9)
    6  = @(8, 8)
    7  = @(6, 9)
    8  = @(7, 9)
    9  = @(6, 6)
    10 = @(9, 9)
    11 = @(8, 8)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]   # column I
    $ws.Cells.Item($row, 10).Value = $pair[1]  # column J
}
